$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 230, pushing existing rows 230-334 down to 231-335.
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new record.
$ws.Range("A230").Value = 3
$ws.Range("B230").Value = "Femacal de La Calera"
$ws.Range("C230").Value = "Coquimbo"
$ws.Range("D230").Value = 44489
$ws.Range("E230").Value = 5
$ws.Range("F230").Value = 100112045
$ws.Range("G230").Value = "Zapallo"
$ws.Range("H230").Value = "Camote"
$ws.Range("I230").Value = "1a (guarda)"
$ws.Range("J230").Value = 220
$ws.Range("K230").Value = 550
$ws.Range("L230").Value = 600
$ws.Range("M230").Value = 577
$ws.Range("N230").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O230").Value = "Provincia de Talca"
$ws.Range("P230").Value = 577
$ws.Range("Q230").Value = 1
$ws.Range("R230").Value = "Hortaliza"
